# Update "想去人数" (F column) counts on the 展览 and 全部类型 sheets
# to match the newly scraped figures (gh-pages output at 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 5497
$ws.Range("F4").Value = 36
$ws.Range("F6").Value = 28
$ws.Range("F7").Value = 646
$ws.Range("F8").Value = 624
$ws.Range("F9").Value = 1070
$ws.Range("F10").Value = 227
$ws.Range("F11").Value = 1526
$ws.Range("F12").Value = 4952
$ws.Range("F13").Value = 449
$ws.Range("F14").Value = 220
$ws.Range("F15").Value = 195
$ws.Range("F16").Value = 5
$ws.Range("F17").Value = 106
$ws.Range("F18").Value = 4268
$ws.Range("F19").Value = 195
$ws.Range("F20").Value = 1139
$ws.Range("F21").Value = 114
$ws.Range("F22").Value = 51
$ws.Range("F23").Value = 207
$ws.Range("F24").Value = 49
$ws.Range("F25").Value = 152
$ws.Range("F27").Value = 145
$ws.Range("F28").Value = 79
$ws.Range("F29").Value = 339
$ws.Range("F30").Value = 39
$ws.Range("F31").Value = 65
$ws.Range("F32").Value = 28
$ws.Range("F33").Value = 40
$ws.Range("F34").Value = 41

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 5497
$ws.Range("F5").Value = 36
$ws.Range("F7").Value = 28
$ws.Range("F8").Value = 646
$ws.Range("F9").Value = 624
$ws.Range("F10").Value = 1070
$ws.Range("F11").Value = 227
$ws.Range("F12").Value = 1526
$ws.Range("F13").Value = 4952
$ws.Range("F14").Value = 449
$ws.Range("F15").Value = 220
$ws.Range("F16").Value = 195
$ws.Range("F17").Value = 5
$ws.Range("F18").Value = 106
$ws.Range("F19").Value = 4268
$ws.Range("F20").Value = 195
$ws.Range("F21").Value = 1139
$ws.Range("F22").Value = 114
$ws.Range("F23").Value = 51
$ws.Range("F24").Value = 207
$ws.Range("F25").Value = 49
$ws.Range("F26").Value = 152
$ws.Range("F27").Value = 55
$ws.Range("F28").Value = 145
$ws.Range("F29").Value = 79
$ws.Range("F30").Value = 339
$ws.Range("F31").Value = 39
$ws.Range("F32").Value = 65
$ws.Range("F33").Value = 28
$ws.Range("F34").Value = 40
$ws.Range("F35").Value = 41

